$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Shift the tail of the metadata table (Description/Purpose/Copyright/Immutable,
# rows 12-15) down by one row to make room for a new "Jurisdiction" row,
# copying bottom-up so we never clobber a row before it has been read.
for ($r = 15; $r -ge 12; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value2 = $ws.Cells.Item($r, 2).Value2
}

# The newly created row 16 needs the same formatting (style index 2 / border +
# wrap) as the rest of the table body - clone it from the row right above.
$ws.Range("A15:B15").Copy() | Out-Null
$ws.Range("A16:B16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Update the changed metadata values ---

# Version bump
$ws.Range("B3").Value2 = "0.1.7"

# Status moved from active to draft
$ws.Range("B6").Value2 = "draft"

# New publish date
$ws.Range("B8").Value2 = "2024-11-22T12:33:30-06:00"

# First Contact row now shows the org contact detail (row 10)
$ws.Range("B10").Value2 = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Second Contact row now shows the individual contact (row 11)
$ws.Range("B11").Value2 = "Bob Milius (bmilius@nmdp.org)"

# New Jurisdiction row (row 12), value left blank
$ws.Range("A12").Value2 = "Jurisdiction"
$ws.Range("B12").Value2 = ""

"done"
